$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: change JobTitle in B5 from "Sandler" to "President"
$ws.Range("B5").Value = "President"

# Row 6: change FirstName in A6 from "Nancy" to the new shared string "Madam"
$ws.Range("A6").Value = "Madam"
